$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure columns -------------------------------------------------
# Original header row: A=canal B=nombres C=apellidos D=correo
#                       E=establecimientoProveniente F=estado G=captador
#
# Target header row:    A=canal B=cedula C=nombres D=apellidos E=celular
#                       F=casa G=correo H=establecimientoProveniente I=captador

# 1) Insert a new column before "nombres" for "cedula"
$ws.Columns("B").Insert()
$ws.Range("B1").Value = "cedula"

# 2) Remove the old "estado" column entirely (now shifted to column G)
$ws.Columns("G").Delete()

# 3) Insert two new columns before "correo" (now column D) for
#    "celular" and "casa"
$ws.Columns("E").Insert()
$ws.Columns("F").Insert()
$ws.Range("E1").Value = "celular"
$ws.Range("F1").Value = "casa"

# --- Add the new data row -------------------------------------------------
$ws.Range("A2").Value = "pega"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "092"

$ws.Range("C2").Value = "jda"
$ws.Range("D2").Value = "jas"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "01923"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "0923"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "asdf"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "asf"

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "asdf"

# --- Misc view state -------------------------------------------------------
$ws.Range("B3").Select() | Out-Null
